# Auto-generated: refresh market-data columns (H:N) across all 8 Leve sheets
# per scheduled-runner diff. Deletions (old cell removed entirely) are applied
# by clearing the cell value to "" which drops the <c> node on save.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 3907.2856
$ws.Cells.Item(2, 9).Value = 4709.091
$ws.Cells.Item(2, 10).Value = 967.3333
$ws.Cells.Item(2, 11).Value = 4709.091
$ws.Cells.Item(2, 12).Value = 967.3333
$ws.Cells.Item(2, 13).Value = -4596.091
$ws.Cells.Item(2, 14).Value = -1193.3333
$ws.Cells.Item(19, 8).Value = 1733
$ws.Cells.Item(19, 10).Value = 1541.25
$ws.Cells.Item(19, 12).Value = 1541.25
$ws.Cells.Item(19, 14).Value = -1891.25
$ws.Cells.Item(32, 8).Value = 950
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 10).Value = 950
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 12).Value = 950
$ws.Cells.Item(32, 13).Value = ""
$ws.Cells.Item(32, 14).Value = -1602
$ws.Cells.Item(40, 8).Value = 2026.6316
$ws.Cells.Item(40, 10).Value = 2155
$ws.Cells.Item(40, 12).Value = 2155
$ws.Cells.Item(40, 14).Value = -2505
$ws.Cells.Item(64, 8).Value = 39835.258
$ws.Cells.Item(64, 10).Value = 2918.5264
$ws.Cells.Item(64, 12).Value = 2918.5264
$ws.Cells.Item(64, 14).Value = -3414.5264
$ws.Cells.Item(67, 8).Value = 39835.258
$ws.Cells.Item(67, 10).Value = 2918.5264
$ws.Cells.Item(67, 12).Value = 2918.5264
$ws.Cells.Item(67, 14).Value = -4634.526400000001
$ws.Cells.Item(74, 8).Value = 3088.9565
$ws.Cells.Item(74, 9).Value = 3071.7334
$ws.Cells.Item(74, 11).Value = 3071.7334
$ws.Cells.Item(74, 13).Value = -2135.7334
$ws.Cells.Item(77, 8).Value = 3088.9565
$ws.Cells.Item(77, 9).Value = 3071.7334
$ws.Cells.Item(77, 11).Value = 15358.667
$ws.Cells.Item(77, 13).Value = -10678.667
$ws.Cells.Item(111, 8).Value = 6694
$ws.Cells.Item(111, 9).Value = 8506
$ws.Cells.Item(111, 10).Value = 2293.4285
$ws.Cells.Item(111, 11).Value = 25518
$ws.Cells.Item(111, 12).Value = 6880.2855
$ws.Cells.Item(111, 13).Value = -22451
$ws.Cells.Item(111, 14).Value = -13014.2855
$ws.Cells.Item(121, 8).Value = 1197.75
$ws.Cells.Item(121, 10).Value = 1083.1428
$ws.Cells.Item(121, 12).Value = 3249.4284
$ws.Cells.Item(121, 14).Value = -6743.428400000001
$ws.Cells.Item(132, 8).Value = 5561056
$ws.Cells.Item(132, 9).Value = 6103042
$ws.Cells.Item(132, 10).Value = 5702.75
$ws.Cells.Item(132, 11).Value = 18309126
$ws.Cells.Item(132, 12).Value = 17108.25
$ws.Cells.Item(132, 13).Value = -18306596
$ws.Cells.Item(132, 14).Value = -22168.25
$ws.Cells.Item(138, 8).Value = 1812.2
$ws.Cells.Item(138, 9).Value = 1976.8572
$ws.Cells.Item(138, 10).Value = 1723.5385
$ws.Cells.Item(138, 11).Value = 5930.571599999999
$ws.Cells.Item(138, 12).Value = 5170.6155
$ws.Cells.Item(138, 13).Value = -790.5715999999993
$ws.Cells.Item(138, 14).Value = -15450.6155

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1694.1842
$ws.Cells.Item(61, 9).Value = 1313.6072
$ws.Cells.Item(61, 10).Value = 2759.8
$ws.Cells.Item(61, 11).Value = 1313.6072
$ws.Cells.Item(61, 12).Value = 2759.8
$ws.Cells.Item(61, 13).Value = -1101.6072
$ws.Cells.Item(61, 14).Value = -3183.8
$ws.Cells.Item(74, 8).Value = 930.1053000000001
$ws.Cells.Item(74, 9).Value = 838
$ws.Cells.Item(74, 10).Value = 1188
$ws.Cells.Item(74, 11).Value = 838
$ws.Cells.Item(74, 12).Value = 1188
$ws.Cells.Item(74, 13).Value = 36
$ws.Cells.Item(74, 14).Value = -2936
$ws.Cells.Item(77, 8).Value = 930.1053000000001
$ws.Cells.Item(77, 9).Value = 838
$ws.Cells.Item(77, 10).Value = 1188
$ws.Cells.Item(77, 11).Value = 4190
$ws.Cells.Item(77, 12).Value = 5940
$ws.Cells.Item(77, 13).Value = 178
$ws.Cells.Item(77, 14).Value = -14676
$ws.Cells.Item(106, 8).Value = 44495
$ws.Cells.Item(106, 10).Value = 44495
$ws.Cells.Item(106, 12).Value = 44495
$ws.Cells.Item(106, 14).Value = -47019
$ws.Cells.Item(122, 8).Value = 1225.4333
$ws.Cells.Item(122, 9).Value = 1016.8182
$ws.Cells.Item(122, 11).Value = 3050.4546
$ws.Cells.Item(122, 13).Value = -600.4546
$ws.Cells.Item(132, 8).Value = 3087.46
$ws.Cells.Item(132, 9).Value = 2895.375
$ws.Cells.Item(132, 11).Value = 8686.125
$ws.Cells.Item(132, 13).Value = -6156.125
$ws.Cells.Item(136, 8).Value = 1694.1842
$ws.Cells.Item(136, 9).Value = 1313.6072
$ws.Cells.Item(136, 10).Value = 2759.8
$ws.Cells.Item(136, 11).Value = 3940.8216
$ws.Cells.Item(136, 12).Value = 8279.400000000001
$ws.Cells.Item(136, 13).Value = -1390.8216
$ws.Cells.Item(136, 14).Value = -13379.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 4004000
$ws.Cells.Item(4, 10).Value = 1006000
$ws.Cells.Item(4, 12).Value = 1006000
$ws.Cells.Item(4, 14).Value = -1006224
$ws.Cells.Item(29, 8).Value = 19500
$ws.Cells.Item(29, 10).Value = 19500
$ws.Cells.Item(29, 12).Value = 19500
$ws.Cells.Item(29, 14).Value = -20086
$ws.Cells.Item(31, 8).Value = 39905.156
$ws.Cells.Item(31, 9).Value = 796.53845
$ws.Cells.Item(31, 10).Value = 60241.64
$ws.Cells.Item(31, 11).Value = 796.53845
$ws.Cells.Item(31, 12).Value = 60241.64
$ws.Cells.Item(31, 13).Value = -501.53845
$ws.Cells.Item(31, 14).Value = -60831.64
$ws.Cells.Item(34, 8).Value = 39905.156
$ws.Cells.Item(34, 9).Value = 796.53845
$ws.Cells.Item(34, 10).Value = 60241.64
$ws.Cells.Item(34, 11).Value = 796.53845
$ws.Cells.Item(34, 12).Value = 60241.64
$ws.Cells.Item(34, 13).Value = -594.53845
$ws.Cells.Item(34, 14).Value = -60645.64
$ws.Cells.Item(86, 8).Value = 3844.1667
$ws.Cells.Item(86, 9).Value = 3450
$ws.Cells.Item(86, 10).Value = 3956.7856
$ws.Cells.Item(86, 11).Value = 3450
$ws.Cells.Item(86, 12).Value = 3956.7856
$ws.Cells.Item(86, 13).Value = -2327
$ws.Cells.Item(86, 14).Value = -6202.7856
$ws.Cells.Item(89, 8).Value = 3844.1667
$ws.Cells.Item(89, 9).Value = 3450
$ws.Cells.Item(89, 10).Value = 3956.7856
$ws.Cells.Item(89, 11).Value = 17250
$ws.Cells.Item(89, 12).Value = 19783.928
$ws.Cells.Item(89, 13).Value = -11634
$ws.Cells.Item(89, 14).Value = -31015.928
$ws.Cells.Item(132, 8).Value = 3664.4167
$ws.Cells.Item(132, 9).Value = 3813.0527
$ws.Cells.Item(132, 10).Value = 3099.6
$ws.Cells.Item(132, 11).Value = 11439.1581
$ws.Cells.Item(132, 12).Value = 9298.799999999999
$ws.Cells.Item(132, 13).Value = -8909.158100000001
$ws.Cells.Item(132, 14).Value = -14358.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1449.2307
$ws.Cells.Item(4, 9).Value = 210
$ws.Cells.Item(4, 10).Value = 2000
$ws.Cells.Item(4, 11).Value = 630
$ws.Cells.Item(4, 12).Value = 6000
$ws.Cells.Item(4, 13).Value = -518
$ws.Cells.Item(4, 14).Value = -6224
$ws.Cells.Item(8, 8).Value = 221.5
$ws.Cells.Item(8, 9).Value = 221.5
$ws.Cells.Item(8, 11).Value = 664.5
$ws.Cells.Item(8, 13).Value = -525.5
$ws.Cells.Item(38, 8).Value = 96
$ws.Cells.Item(38, 9).Value = 80
$ws.Cells.Item(38, 10).Value = 98.666664
$ws.Cells.Item(38, 11).Value = 240
$ws.Cells.Item(38, 12).Value = 295.999992
$ws.Cells.Item(38, 13).Value = 107
$ws.Cells.Item(38, 14).Value = -989.999992
$ws.Cells.Item(41, 8).Value = 2267.6667
$ws.Cells.Item(41, 10).Value = 3151.5
$ws.Cells.Item(41, 12).Value = 9454.5
$ws.Cells.Item(41, 14).Value = -10130.5
$ws.Cells.Item(42, 8).Value = 2602.5
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 2602.5
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 12).Value = 7807.5
$ws.Cells.Item(42, 13).Value = ""
$ws.Cells.Item(42, 14).Value = -8875.5
$ws.Cells.Item(43, 8).Value = 3003
$ws.Cells.Item(43, 10).Value = 3003
$ws.Cells.Item(43, 12).Value = 9009
$ws.Cells.Item(43, 14).Value = -9237
$ws.Cells.Item(131, 8).Value = 882.89
$ws.Cells.Item(131, 10).Value = 898.85565
$ws.Cells.Item(131, 12).Value = 2696.56695
$ws.Cells.Item(131, 14).Value = -12776.56695

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2626.2
$ws.Cells.Item(102, 9).Value = 2314.923
$ws.Cells.Item(102, 10).Value = 3204.2856
$ws.Cells.Item(102, 11).Value = 2314.923
$ws.Cells.Item(102, 12).Value = 3204.2856
$ws.Cells.Item(102, 13).Value = -692.9229999999998
$ws.Cells.Item(102, 14).Value = -6448.2856
$ws.Cells.Item(132, 8).Value = 3919.5715
$ws.Cells.Item(132, 9).Value = 2080
$ws.Cells.Item(132, 10).Value = 5299.25
$ws.Cells.Item(132, 11).Value = 6240
$ws.Cells.Item(132, 12).Value = 15897.75
$ws.Cells.Item(132, 13).Value = -3710
$ws.Cells.Item(132, 14).Value = -20957.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(123, 8).Value = 23442.5
$ws.Cells.Item(123, 9).Value = 14390
$ws.Cells.Item(123, 10).Value = 32495
$ws.Cells.Item(123, 11).Value = 14390
$ws.Cells.Item(123, 12).Value = 32495
$ws.Cells.Item(123, 13).Value = -9490
$ws.Cells.Item(123, 14).Value = -42295
$ws.Cells.Item(136, 8).Value = 1500.1177
$ws.Cells.Item(136, 9).Value = 1251.7307
$ws.Cells.Item(136, 10).Value = 2307.375
$ws.Cells.Item(136, 11).Value = 3755.1921
$ws.Cells.Item(136, 12).Value = 6922.125
$ws.Cells.Item(136, 13).Value = -1205.1921
$ws.Cells.Item(136, 14).Value = -12022.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(56, 8).Value = 16666.666
$ws.Cells.Item(56, 9).Value = 3500
$ws.Cells.Item(56, 10).Value = 43000
$ws.Cells.Item(56, 11).Value = 3500
$ws.Cells.Item(56, 12).Value = 43000
$ws.Cells.Item(56, 13).Value = -2786
$ws.Cells.Item(56, 14).Value = -44428
$ws.Cells.Item(107, 8).Value = 72008
$ws.Cells.Item(107, 9).Value = 456.8889
$ws.Cells.Item(107, 10).Value = 200800
$ws.Cells.Item(107, 11).Value = 1370.6667
$ws.Cells.Item(107, 12).Value = 602400
$ws.Cells.Item(107, 13).Value = 549.3333
$ws.Cells.Item(107, 14).Value = -606240
$ws.Cells.Item(126, 8).Value = 2048.8
$ws.Cells.Item(126, 9).Value = 2038.8182
$ws.Cells.Item(126, 10).Value = 2076.25
$ws.Cells.Item(126, 11).Value = 6116.4546
$ws.Cells.Item(126, 12).Value = 6228.75
$ws.Cells.Item(126, 13).Value = -3646.4546
$ws.Cells.Item(126, 14).Value = -11168.75
